# "try overflow for mobile"
# - Fill in the missing column D (price) values for the "Tourist" rows that
#   were left blank (rows 97-113 and 609-625 on Sheet1).
# - Scroll the view down (topLeftCell / selection) as if someone had scrolled
#   to inspect that area on a mobile-sized window.
# - Turn on AutoFilter for the used range, which also registers the sheet's
#   hidden _FilterDatabase defined name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D = 4
$updates = @(
    @{Row=97; Value=6.8},
    @{Row=98; Value=6.8},
    @{Row=99; Value=7.9},
    @{Row=100; Value=5.9},
    @{Row=101; Value=7.9},
    @{Row=102; Value=6.65},
    @{Row=103; Value=5.95},
    @{Row=104; Value=5.65},
    @{Row=105; Value=7.3},
    @{Row=106; Value=6.4},
    @{Row=107; Value=4.5},
    @{Row=108; Value=5.45},
    @{Row=109; Value=5.25},
    @{Row=110; Value=4.5},
    @{Row=111; Value=4.5},
    @{Row=112; Value=7.1},
    @{Row=113; Value=5.5},
    @{Row=609; Value=6.8},
    @{Row=610; Value=6.8},
    @{Row=611; Value=8.5},
    @{Row=612; Value=5.9},
    @{Row=613; Value=9.4},
    @{Row=614; Value=8.0},
    @{Row=615; Value=7.15},
    @{Row=616; Value=6.85},
    @{Row=617; Value=8.3},
    @{Row=618; Value=7.7},
    @{Row=619; Value=5.3},
    @{Row=620; Value=6.55},
    @{Row=621; Value=6.55},
    @{Row=622; Value=5.5},
    @{Row=623; Value=5.9},
    @{Row=624; Value=8.3},
    @{Row=625; Value=6.75}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value = $u.Value
}

# Turn on the worksheet AutoFilter over the whole used range (also creates
# the hidden _FilterDatabase defined name referenced in workbook.xml).
$ws.Range("A1:D641").AutoFilter()

# Scroll the sheetView down and move the active selection, as if checking
# the layout further down the list (e.g. on a small/mobile viewport).
$ws.Range("A500").Select()
$excel.ActiveWindow.ScrollRow = 500
$ws.Range("C509").Select()

$wb.Save()
